# Apply the "Discount Rate" workbook update (WRI China / Hong Kong EPS v2.0.0 bring-up).
#
# Summary of the functional edit:
#  - "About" sheet: rewrite the Notes section. The old note (3 lines, starting at
#    row 10) explaining the source of the 3% discount rate is kept but moved down
#    to rows 16-18, and a brand new note (5 lines, rows 10-14) is inserted above it
#    explaining what "Annual Perc" actually means.
#  - "DR" sheet: the header label in B1 gets a parenthetical unit suffix
#    ("Annual Perc" -> "Annual Perc (dimensionless)"), the cell wraps its text, and
#    row 1 is made taller (30pt) to fit the now two-line label.
#  - Final on-screen selections are left on A16:A18 on "About" and B1 on "DR",
#    with "About" remaining the active sheet/tab.

$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsDR = $wb.Worksheets.Item("DR")

# ---- About sheet: replace + extend the Notes section ----------------------

# New explanatory note about what "Annual Perc" means (rows 10-14).
$wsAbout.Range("A10").Value = "This is the annual percentage rate by which future savings (e.g. fuel cost savings) are discounted when"
$wsAbout.Range("A11").Value = "making price-driven purchasing decisions in the current year. The value used should be one that is"
$wsAbout.Range("A12").Value = "reasonable for people who are looking to buy fuel-consuming capital equipment, such as industrial"
$wsAbout.Range("A13").Value = "equipment or building components. The model works in real dollars, so this rate should be the growth"
$wsAbout.Range("A14").Value = "in real value, not the growth in nominal value plus real value."

# Row 15 stays blank (spacer), matching the blank-row style already used after row 8.

# The original 3-line note about the 3% discount rate source moves down to 16-18.
$wsAbout.Range("A16").Value = "We choose to use a 3% discount rate here, for consistency with the 3% rate used for the central estimate"
$wsAbout.Range("A17").Value = "of Social Cost of Carbon (in the SCoC variable), as well as the discount rate built into the health"
$wsAbout.Range("A18").Value = "damages values in the SCoHIbP Social Cost of Health Impacts by Pollutant variable."

# ---- DR sheet: relabel + reformat the header cell --------------------------

$wsDR.Range("B1").Value = "Annual Perc (dimensionless)"
$wsDR.Range("B1").WrapText = $true
$wsDR.Rows.Item(1).RowHeight = 30

# ---- Final selections (About stays the active/visible sheet) ---------------

$wsDR.Range("B1").Select()
$wsAbout.Range("A16:A18").Select()
